# Daily attendance processing - 2026-01-23 11:10:06
#
# The "Recorded By" column (G) lists the users who recorded attendance for
# a session. Cells that previously read "System, dnasr281@gmail.com" need
# to have the author order flipped to "dnasr281@gmail.com, System".
#
# Every cell in the used range whose value is the exact original string is
# updated - this mirrors the source edit, which touched only those G-column
# cells (rows where System was recorded before the human user) and left
# every other "Recorded By" value (e.g. "dnasr281@gmail.com" alone, or
# "System" alone) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$used.Replace("System, dnasr281@gmail.com", "dnasr281@gmail.com, System")
